$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Granada" (Vega Modelo de Temuco) is inserted
# at row 106, pushing the existing records (old rows 106-153) down by one row
# (new rows 107-154).
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with the new record's data. The
# descriptive columns (A-C, E-L) match the rest of the Granada block.
$ws.Range("A106").Value = 10
$ws.Range("B106").Value = "Vega Modelo de Temuco"
$ws.Range("C106").Value = "La Araucanía"
$ws.Range("D106").Value = 44784
$ws.Range("E106").Value = 9
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100104
$ws.Range("H106").Value = "Frutos de pepita"
$ws.Range("I106").Value = 100104001
$ws.Range("J106").Value = "Granada"
$ws.Range("K106").Value = "Wonderfull"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 45
$ws.Range("N106").Value = 15000
$ws.Range("O106").Value = 15000
$ws.Range("P106").Value = 15000
$ws.Range("Q106").Value = "$/bandeja 10 kilos granel"
$ws.Range("R106").Value = "Provincia de Limarí"
$ws.Range("S106").Value = 1500
$ws.Range("T106").Value = 10
